$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.000951647758484
$ws.Range("B1").Value = 3.151415824890137
$ws.Range("C1").Value = 3.891799926757812
$ws.Range("D1").Value = 2.029496192932129
$ws.Range("E1").Value = 1.199029445648193
